# Apply updated crypto market data (price + 1h volume change) per the Oct 28 2023 refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '34.189.50'
$ws.Range("E2").Value = '  +1.00%  '
$ws.Range("D3").Value = '1.781.12'
$ws.Range("E3").Value = '  +0.12%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '''226.08'
$ws.Range("E5").Value = '  +0.96%  '
$ws.Range("D6").Value = '''0.547'
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("E7").Value = '  +0.17%  '
$ws.Range("D8").Value = '''31.75'
$ws.Range("E8").Value = '  +0.11%  '
$ws.Range("D9").Value = '''0.292'
$ws.Range("E9").Value = '  +0.80%  '
$ws.Range("D10").Value = '''0.0692'
$ws.Range("E10").Value = '  +2.24%  '
$ws.Range("D11").Value = '''0.0944'
$ws.Range("E11").Value = '  +0.92%  '
$ws.Range("D12").Value = '2.037.27'
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("D13").Value = '''10.96'
$ws.Range("E13").Value = '  -1.79%  '
$ws.Range("D14").Value = '1.779.38'
$ws.Range("E14").Value = '  +0.31%  '
$ws.Range("D15").Value = '34.157.29'
$ws.Range("E15").Value = '  +0.82%  '
$ws.Range("E16").Value = '  +2.36%  '
$ws.Range("D17").Value = '''4.18'
$ws.Range("E17").Value = '  +1.06%  '
$ws.Range("D18").Value = '''67.89'
$ws.Range("E18").Value = '  +1.76%  '
$ws.Range("D19").Value = '0.0₃0802'
$ws.Range("E19").Value = '  +4.13%  '
$ws.Range("D20").Value = '''246.05'
$ws.Range("E20").Value = '  +2.97%  '
$ws.Range("D21").Value = '''10.97'
$ws.Range("E21").Value = '  +3.93%  '
$ws.Range("E22").Value = '  +0.17%  '
$ws.Range("D23").Value = '''4.09'
$ws.Range("E23").Value = '  +1.93%  '
$ws.Range("E24").Value = '  -1.27%  '
$ws.Range("D25").Value = '''162.31'
$ws.Range("D26").Value = '''7.19'
$ws.Range("E26").Value = '  +2.25%  '
$ws.Range("D27").Value = '''16.30'
$ws.Range("E27").Value = '  +1.36%  '
$ws.Range("E28").Value = '  +1.80%  '
$ws.Range("D29").Value = '''1.00'
$ws.Range("E29").Value = '  +0.18%  '
$ws.Range("E30").Value = '  +1.02%  '
$ws.Range("D31").Value = '''0.0521'
$ws.Range("E31").Value = '  +2.35%  '
$ws.Range("D32").Value = '''3.74'
$ws.Range("E32").Value = '  +4.35%  '
$ws.Range("E33").Value = '  +5.53%  '
$ws.Range("E34").Value = '  -1.02%  '
$ws.Range("D35").Value = '1.441.55'
$ws.Range("E35").Value = '  +3.70%  '
$ws.Range("D36").Value = '''0.657'
$ws.Range("E36").Value = '  +3.76%  '
$ws.Range("D37").Value = '''2.40'
$ws.Range("E37").Value = '  +6.77%  '
$ws.Range("E38").Value = '  +3.00%  '
$ws.Range("E39").Value = '  +0.29%  '
$ws.Range("D40").Value = '''80.15'
$ws.Range("E40").Value = '  +2.42%  '
$ws.Range("E41").Value = '  -0.12%  '
$ws.Range("D42").Value = '''0.923'
$ws.Range("E42").Value = '  +1.44%  '
$ws.Range("E43").Value = '  +0.70%  '
$ws.Range("D44").Value = '''13.53'
$ws.Range("E44").Value = '  +0.47%  '
$ws.Range("B45").Value = 'FraxShare'
$ws.Range("C45").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D45").Value = '''6.08'
$ws.Range("E45").Value = '  +3.74%  '
$ws.Range("B46").Value = 'Kaspa'
$ws.Range("C46").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D46").Value = '''0.0510'
$ws.Range("E46").Value = '  +0.24%  '
$ws.Range("E47").Value = '  +0.11%  '
$ws.Range("E48").Value = '  -2.92%  '
$ws.Range("D49").Value = '1.940.22'
$ws.Range("E49").Value = '  +0.14%  '
$ws.Range("D50").Value = '''104.43'
$ws.Range("E50").Value = '  -1.54%  '
$ws.Range("E51").Value = '  +0.25%  '
